# Updated cryptos list with refreshed Price (D) and Volume(1h) (E) values.
# Values are written through a TEXT-producing formula (=T("...")) and then
# "pasted as values" (xlPasteValues = -4163) so the resulting cells hold
# plain text strings (matching the source inline-string cells) instead of
# being auto-coerced into numbers by Excel, and without touching any
# cell styles/number formats.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$c = $ws.Range('D2')
$c.Formula = '=T("62.847.53")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E2')
$c.Formula = '=T("  +1.24%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D3')
$c.Formula = '=T("2.436.54")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E3')
$c.Formula = '=T("  +0.82%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E4')
$c.Formula = '=T("  +0.08%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D5')
$c.Formula = '=T("570.51")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E5')
$c.Formula = '=T("  +1.28%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D6')
$c.Formula = '=T("146.42")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E6')
$c.Formula = '=T("  +2.58%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D7')
$c.Formula = '=T("1.00")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E7')
$c.Formula = '=T("  -0.01%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E8')
$c.Formula = '=T("  +0.78%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D9')
$c.Formula = '=T("0.111")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E9')
$c.Formula = '=T("  +1.78%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E10')
$c.Formula = '=T("  +0.49%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D11')
$c.Formula = '=T("5.31")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E11')
$c.Formula = '=T("  +1.58%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D12')
$c.Formula = '=T("0.357")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E12')
$c.Formula = '=T("  +1.88%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D13')
$c.Formula = '=T("26.92")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E13')
$c.Formula = '=T("  +4.73%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D14')
$c.Formula = '=T("0.0000181")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E14')
$c.Formula = '=T("  +3.88%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D15')
$c.Formula = '=T("2.879.06")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E15')
$c.Formula = '=T("  +0.90%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D16')
$c.Formula = '=T("62.687.88")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E16')
$c.Formula = '=T("  +1.10%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D17')
$c.Formula = '=T("2.440.45")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E17')
$c.Formula = '=T("  +1.07%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D18')
$c.Formula = '=T("11.28")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E18')
$c.Formula = '=T("  +0.18%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D19')
$c.Formula = '=T("7.07")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E19')
$c.Formula = '=T("  +3.60%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D20')
$c.Formula = '=T("324.99")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E20')
$c.Formula = '=T("  +0.97%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E21')
$c.Formula = '=T("  +1.22%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E22')
$c.Formula = '=T("  +0.18%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D23')
$c.Formula = '=T("1.84")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E23')
$c.Formula = '=T("  +4.91%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D24')
$c.Formula = '=T("67.18")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E24')
$c.Formula = '=T("  +1.51%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D25')
$c.Formula = '=T("634.38")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E25')
$c.Formula = '=T("  +11.96%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D26')
$c.Formula = '=T("8.69")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E26')
$c.Formula = '=T("  +0.88%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D27')
$c.Formula = '=T("0.0000102")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E27')
$c.Formula = '=T("  +9.42%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D28')
$c.Formula = '=T("2.557.77")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E28')
$c.Formula = '=T("  +0.98%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E29')
$c.Formula = '=T("  +3.67%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E30')
$c.Formula = '=T("  -0.11%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E31')
$c.Formula = '=T("  +4.64%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D32')
$c.Formula = '=T("0.142")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E32')
$c.Formula = '=T("  -3.66%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E33')
$c.Formula = '=T("  +0.97%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E34')
$c.Formula = '=T("  +0.24%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D35')
$c.Formula = '=T("4.96")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E35')
$c.Formula = '=T("  +3.76%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E36')
$c.Formula = '=T("  -0.01%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E37')
$c.Formula = '=T("  +0.78%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D38')
$c.Formula = '=T("5.44")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E38')
$c.Formula = '=T("  -0.05%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D39')
$c.Formula = '=T("18.76")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E39')
$c.Formula = '=T("  +1.09%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E40')
$c.Formula = '=T("  +1.89%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D41')
$c.Formula = '=T("148.26")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E41')
$c.Formula = '=T("  -3.00%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D42')
$c.Formula = '=T("2.57")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E42')
$c.Formula = '=T("  +14.20%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E43')
$c.Formula = '=T("  +0.29%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D44')
$c.Formula = '=T("150.08")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E44')
$c.Formula = '=T("  +0.88%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D45')
$c.Formula = '=T("3.69")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E45')
$c.Formula = '=T("  +2.19%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E46')
$c.Formula = '=T("  +1.82%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('D47')
$c.Formula = '=T("20.77")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E47')
$c.Formula = '=T("  +4.34%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E48')
$c.Formula = '=T("  +1.84%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E50')
$c.Formula = '=T("  +0.44%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range('E51')
$c.Formula = '=T("  +4.34%  ")'
$c.Copy()
$c.PasteSpecial(-4163)
$excel.CutCopyMode = 0
